# Add two new daily-log rows (13 and 14) to Sheet1, mirroring the existing
# table layout: column A = day number, column B = time range, column C = content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13
$ws.Range("A13").Value = 29
$ws.Range("B13").Value = "4：42—6：54"
$ws.Range("C13").Value = "dowhile n的阶乘 阶乘和 有序数组寻找n（二分法）"

# Row 14
$ws.Range("A14").Value = 30
$ws.Range("B14").Value = "3：00-7：20"
$ws.Range("C14").Value = "多字符两端移动 输入密码 猜数字游戏"

# Update the view selection to match the new last cell (mirrors the diff's
# sheetView/selection change from C12 to C14).
$ws.Range("C14").Select()
